# Se agregan validaciones en form retenciones - validacion codigo respuesta - gestion de cdr
# Adds a new "Tiene Igv" column (M) with SI/NO values to the item_sets sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + values in column M
$ws.Range("M1").Value = "Tiene Igv"
$ws.Range("M2").Value = "SI"
$ws.Range("M3").Value = "NO"

# Give the new column its own (3rd) cell format, matching the workbook's
# existing pattern of cells carrying an explicit, distinct style index.
$ws.Range("M1:M3").Locked = $true

# Move the selection the way the saved file shows it
$ws.Range("A3").Select()
